# Add a new row for "climate_change_factor_gnrl_hydropower_availability"
# above the existing "elasticity_gnrl_rate_occupancy_to_gdppc" row (row 4),
# pushing all subsequent General/Economy variable rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4 (shifts rows 4:11 down to 5:12)
$ws.Rows("4:4").Insert()

# Populate the new row 4 with the climate change factor variable data
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

# Columns J4:AS4 are all set to 1
$ws.Range("J4:AS4").Value = 1

Write-Output "done"
